$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 2
$ws.Range("H2").Value = 17477.223
$ws.Range("I2").Value = 1933.3334
$ws.Range("K2").Value = 1933.3334
$ws.Range("M2").Value = -1820.3334

# Row 41
$ws.Range("H41").Value = 405.33334
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# Row 57
$ws.Range("H57").Value = 81499
$ws.Range("I57").Value = 67998
$ws.Range("J57").Value = 95000
$ws.Range("K57").Value = 203994
$ws.Range("L57").Value = 285000
$ws.Range("N57").Value = -285998
$ws.Range("M57").Value = -203495

# Row 62
$ws.Range("H62").Value = 5899.5
$ws.Range("I62").Value = 5899.5
$ws.Range("J62").Value = 5899.5
$ws.Range("K62").Value = 5899.5
$ws.Range("L62").Value = 5899.5
$ws.Range("M62").Value = -5275.5
$ws.Range("N62").Value = -7147.5

# Row 65
$ws.Range("H65").Value = 5899.5
$ws.Range("I65").Value = 5899.5
$ws.Range("J65").Value = 5899.5
$ws.Range("K65").Value = 29497.5
$ws.Range("L65").Value = 29497.5
$ws.Range("M65").Value = -26377.5
$ws.Range("N65").Value = -35737.5

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 132
$ws.Range("H132").Value = 962.3077
$ws.Range("I132").Value = 1080.6818
$ws.Range("K132").Value = 3242.0454
$ws.Range("M132").Value = -712.0454

# Row 137
$ws.Range("H137").Value = 2145.3333
$ws.Range("I137").Value = 1468.25
$ws.Range("K137").Value = 4404.75
$ws.Range("M137").Value = -1854.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1153
$ws.Range("I2").Value = 1153
$ws.Range("K2").Value = 1153
$ws.Range("M2").Value = -1040

# Row 32
$ws.Range("H32").Value = 4349.7393
$ws.Range("I32").Value = 3638.4092
$ws.Range("J32").Value = 19999
$ws.Range("K32").Value = 3638.4092
$ws.Range("L32").Value = 19999
$ws.Range("M32").Value = -3351.4092
$ws.Range("N32").Value = -20573

# Row 61
$ws.Range("H61").Value = 3998.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 3998.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 3998.5
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -4422.5

# Row 74
$ws.Range("H74").Value = 4998.5
$ws.Range("I74").Value = 3997
$ws.Range("K74").Value = 3997
$ws.Range("M74").Value = -3123

# Row 77
$ws.Range("H77").Value = 4998.5
$ws.Range("I77").Value = 3997
$ws.Range("K77").Value = 19985
$ws.Range("M77").Value = -15617

# Row 102
$ws.Range("H102").Value = 1214.1428
$ws.Range("I102").Value = 1254.0834
$ws.Range("K102").Value = 1254.0834
$ws.Range("M102").Value = 367.9166

# Row 116
$ws.Range("H116").Value = 1153
$ws.Range("I116").Value = 1153
$ws.Range("K116").Value = 1153
$ws.Range("M116").Value = 1141

# Row 132
$ws.Range("H132").Value = 3724
$ws.Range("I132").Value = 2977.6667
$ws.Range("K132").Value = 8933.000100000001
$ws.Range("M132").Value = -6403.000100000001

# Row 136
$ws.Range("H136").Value = 3998.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3998.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 11995.5
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -17095.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1153
$ws.Range("I3").Value = 1153
$ws.Range("K3").Value = 1153
$ws.Range("M3").Value = -1039

# Row 134
$ws.Range("H134").Value = 4835.5
$ws.Range("I134").Value = 4835.5
$ws.Range("K134").Value = 14506.5
$ws.Range("M134").Value = -11971.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 3247
$ws.Range("I16").Value = 3373.5
$ws.Range("K16").Value = 3373.5
$ws.Range("M16").Value = -3086.5

# Row 31
$ws.Range("H31").Value = 4280.206
$ws.Range("I31").Value = 1629.409
$ws.Range("K31").Value = 1629.409
$ws.Range("M31").Value = -1334.409

# Row 34
$ws.Range("H34").Value = 4280.206
$ws.Range("I34").Value = 1629.409
$ws.Range("K34").Value = 1629.409
$ws.Range("M34").Value = -1427.409

# Row 113
$ws.Range("H113").Value = 3247
$ws.Range("I113").Value = 3373.5
$ws.Range("K113").Value = 3373.5
$ws.Range("M113").Value = -1203.5

# Row 134
$ws.Range("H134").Value = 2296.8125
$ws.Range("I134").Value = 2053.5715
$ws.Range("K134").Value = 6160.7145
$ws.Range("M134").Value = -3625.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 79672930
$ws.Range("I4").Value = 15000826
$ws.Range("J4").Value = 170213870
$ws.Range("K4").Value = 45002478
$ws.Range("L4").Value = 510641610
$ws.Range("M4").Value = -45002366
$ws.Range("N4").Value = -510641834

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 113
$ws.Range("H113").Value = 1499.5
$ws.Range("I113").Value = 1499.5
$ws.Range("K113").Value = 1499.5
$ws.Range("M113").Value = 670.5

# Row 132
$ws.Range("H132").Value = 5230.8887
$ws.Range("I132").Value = 4847
$ws.Range("K132").Value = 14541
$ws.Range("M132").Value = -12011

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 61
$ws.Range("H61").Value = 1102.2
$ws.Range("I61").Value = 1327.75
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 1327.75
$ws.Range("L61").Value = 200
$ws.Range("M61").Value = -1125.75
$ws.Range("N61").Value = -604

# Row 93
$ws.Range("H93").Value = 924.125
$ws.Range("I93").Value = 920.4286
$ws.Range("J93").Value = 950
$ws.Range("K93").Value = 920.4286
$ws.Range("L93").Value = 950
$ws.Range("M93").Value = 327.5714
$ws.Range("N93").Value = -3446

# Row 100
$ws.Range("H100").Value = 1588.75
$ws.Range("I100").Value = 1588.75
$ws.Range("K100").Value = 1588.75
$ws.Range("M100").Value = -1047.75

# Row 113
$ws.Range("H113").Value = 1102.2
$ws.Range("I113").Value = 1327.75
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 1327.75
$ws.Range("L113").Value = 200
$ws.Range("M113").Value = 842.25
$ws.Range("N113").Value = -4540

# Row 132
$ws.Range("H132").Value = 3134.75
$ws.Range("I132").Value = 1833.3636
$ws.Range("K132").Value = 5500.0908
$ws.Range("M132").Value = -2970.0908

# Row 136
$ws.Range("H136").Value = 31330.053
$ws.Range("J136").Value = 45081.25
$ws.Range("L136").Value = 135243.75
$ws.Range("N136").Value = -140343.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 96
$ws.Range("H96").Value = 3085.889
$ws.Range("I96").Value = 3428
$ws.Range("J96").Value = 1888.5
$ws.Range("K96").Value = 3428
$ws.Range("L96").Value = 1888.5
$ws.Range("M96").Value = -2055
$ws.Range("N96").Value = -4634.5

# Row 113
$ws.Range("H113").Value = 1203
$ws.Range("I113").Value = 1099.6
$ws.Range("K113").Value = 3298.8
$ws.Range("M113").Value = -1128.8

# Row 122
$ws.Range("H122").Value = 1487.4706
$ws.Range("I122").Value = 1487.4706
$ws.Range("K122").Value = 4462.4118
$ws.Range("M122").Value = -2012.4118

# Row 132
$ws.Range("H132").Value = 2450.1765
$ws.Range("I132").Value = 2240.88
$ws.Range("K132").Value = 6722.64
$ws.Range("M132").Value = -4192.64

# Row 136
$ws.Range("H136").Value = 5180.4194
$ws.Range("I136").Value = 5657.16
$ws.Range("J136").Value = 3194
$ws.Range("K136").Value = 16971.48
$ws.Range("L136").Value = 9582
$ws.Range("M136").Value = -14421.48
$ws.Range("N136").Value = -14682
